# The MAPPING sheet's TABLE_NAME column (D) was pointing at "_TEST" shadow
# tables. Point it back at the real production table names.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MAPPING")

$ws.Range("D2").Value = "FRIENDLY"
$ws.Range("D3").Value = "HIGH_DATE_TABLES"
$ws.Range("D4").Value = "SLA_CONFIG"
$ws.Range("D5").Value = "RUN_HISTORY_TABLES"
$ws.Range("D6").Value = "DATA_FEED_CATALOG"
$ws.Range("D7").Value = "DATA_FEED_MAPPING"
